# Update NATMI LR-pair sheet (Fgf8-Fgfr3) with refreshed TPM-based numbers.
# New data has only 4 sending/target combinations (all "from ECs"); the
# previous 3 "from MuSCs" rows go away and are replaced by a 4th "from ECs"
# row targeting a new cluster "Resolving-Mac".

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column order (A..T), one row per inner array.
$data = @(
    @("ECs","Fgf8","Fgfr3","ECs",1,0.3333333333333333,0.1008186666666667,0.302456,1,1,3,1,2.330840333333333,6.992521,0.6715345129768794,0.7003397275969581,0.2349922146195556,2.114929931576,0.6715345129768794,0.7003397275969581),
    @("ECs","Fgf8","Fgfr3","FAPs",1,0.3333333333333333,0.1008186666666667,0.302456,1,1,3,1,0.657666,1.972998,0.1894790521235985,0.1976066831789769,0.066305009232,0.596745083088,0.1894790521235985,0.1976066831789769),
    @("ECs","Fgf8","Fgfr3","MuSCs",1,0.3333333333333333,0.1008186666666667,0.302456,1,1,2,1,0.42828,0.85656,0.123391035029171,0.08578923067523865,0.04317861856,0.25907171136,0.123391035029171,0.08578923067523865),
    @("ECs","Fgf8","Fgfr3","Resolving-Mac",1,0.3333333333333333,0.1008186666666667,0.302456,1,1,1,0.3333333333333333,0.05413033333333334,0.162391,0.01559539987035126,0.01626435854882633,0.005457348032888889,0.049116132296,0.01559539987035126,0.01626435854882633)
)

for ($i = 0; $i -lt $data.Length; $i++) {
    $rowNum = $i + 2
    $rowVals = $data[$i]
    for ($c = 0; $c -lt $rowVals.Length; $c++) {
        $ws.Cells.Item($rowNum, $c + 1).Value = $rowVals[$c]
    }
}

# The refreshed table only has 4 data rows now (rows 2-5); drop the old
# rows 6 and 7 that carried the "from MuSCs" pairs.
$ws.Range("A6:A7").EntireRow.Delete()
